# S-01004 - Avance en los datos de los graficos de performance de agentes
#
# Adds two new logged-hours entries for "Nico" on task "S-01004" / "Desarrollo
# Metricas Agentes" (2010-10-15 and 2010-10-16, 4 hours each) to the "Horas
# insumidas" sheet, fixes two neighbouring dates that were mis-entered, and
# lets the "Earned Value" sheet's SUMIF-driven figures recalculate to reflect
# the additional hours.

$wb = $excel.ActiveWorkbook

$wsEV = $wb.Worksheets.Item("Earned Value")
$wsHI = $wb.Worksheets.Item("Horas insumidas")

# --- Fix the two existing rows that had the wrong date -----------------
$wsHI.Range("B68").Value2 = 40465
$wsHI.Range("B69").Value2 = 40465

# --- New entries: Nico worked on "Desarrollo Metricas Agentes" (S-01004) -
$wsHI.Range("B70").Value2 = 40465
$wsHI.Range("B70").NumberFormat = "d-mmm"
$wsHI.Range("C70").Value2 = "Nico"
$wsHI.Range("D70").Value2 = "Desarrollo Metricas Agentes"
$wsHI.Range("E70").Value2 = "S-01004"
$wsHI.Range("F70").Value2 = 4

$wsHI.Range("B71").Value2 = 40466
$wsHI.Range("B71").NumberFormat = "d-mmm"
$wsHI.Range("C71").Value2 = "Nico"
$wsHI.Range("D71").Value2 = "Desarrollo Metricas Agentes"
$wsHI.Range("E71").Value2 = "S-01004"
$wsHI.Range("F71").Value2 = 4

# Recalculate everything so the SUMIF formulas in "Earned Value" see the
# freshly entered hours.
$excel.CalculateFullRebuild()

# The totals in row 21 aggregate (via SUM) cells that are themselves the
# result of the cross-sheet SUMIF formulas above; force them to be rebuilt
# from scratch so they pick up the new totals.
$wsEV.Range("G21").ClearContents()
$wsEV.Range("H21").ClearContents()
$excel.CalculateFullRebuild()
$wsEV.Range("G21").Formula = "=SUM(G2:G11)"
$wsEV.Range("H21").Formula = "=SUM(H2:H11)"
$excel.CalculateFullRebuild()

# --- Update the view of "Horas insumidas" to where the edit took place --
[void]$wsHI.Activate()
[void]$wsHI.Range("D73").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 60
$win.ScrollColumn = 1
